# Apply "feat: add 2022-Q3 data":
#  - Insert a new worksheet "2022-Q3" (holding per-fund Q3 holdings) right
#    after the "总计" summary sheet and before the existing "2022-Q1" sheet.
#  - Populate the new sheet with the Q3 fund holdings table.
#  - Update the "总计" summary sheet: the top data row becomes the new
#    2022-Q3 totals, and the former 2022-Q1 totals row is pushed down to
#    row 3 (unchanged values).

$wb = $excel.ActiveWorkbook

function Style-HeaderCell($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous
    $cell.Borders.Weight = 2            # xlThin
}

function Set-TextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet, positioned right after "总计".
# ---------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item(1)

$q3Sheet = $wb.Worksheets.Add($null, $summarySheet)
$q3Sheet.Name = "2022-Q3"

# ---------------------------------------------------------------
# 2. Populate the "2022-Q3" sheet with the fund holdings table.
# ---------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q3Sheet.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    Style-HeaderCell $cell
}

$q3Rows = @(
    @(0, "519956", "长信睿进灵活配置混合C", "6.86", "39.13", "2.43", "0.1667", 7),
    @(1, "002210", "创金合信量化多因子股票A", "2.39", "91.71", "1.22", "0.0292", 8),
    @(2, "003865", "创金合信量化多因子股票C", "0.75", "91.71", "1.22", "0.0092", 8),
    @(3, "519957", "长信睿进灵活配置混合A", "0.02", "39.13", "2.43", "0.0005", 7)
)

$r = 2
foreach ($row in $q3Rows) {
    $aCell = $q3Sheet.Cells.Item($r, 1)
    $aCell.Value = $row[0]
    Style-HeaderCell $aCell

    Set-TextCell $q3Sheet $r 2 $row[1]
    Set-TextCell $q3Sheet $r 3 $row[2]
    Set-TextCell $q3Sheet $r 4 $row[3]
    Set-TextCell $q3Sheet $r 5 $row[4]
    Set-TextCell $q3Sheet $r 6 $row[5]
    Set-TextCell $q3Sheet $r 7 $row[6]
    $q3Sheet.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------
# 3. Update the "总计" summary sheet: row 2 becomes the 2022-Q3 totals,
#    and the old 2022-Q1 totals (previously row 2) move down to row 3.
# ---------------------------------------------------------------
$summarySheet.Cells.Item(3, 1).Value = 1
Style-HeaderCell $summarySheet.Cells.Item(3, 1)
Set-TextCell $summarySheet 3 2 "2022-Q1"
$summarySheet.Cells.Item(3, 3).Value = 1
$summarySheet.Cells.Item(3, 4).Value = 0.03

Set-TextCell $summarySheet 2 2 "2022-Q3"
$summarySheet.Cells.Item(2, 3).Value = 4
$summarySheet.Cells.Item(2, 4).Value = 0.21

Write-Host "Applied 2022-Q3 data successfully."
